$p = $ppt.ActivePresentation

# --- Slide 8 ("Peephole Optimization" / redundant-instruction-elimination
#     bullets): the "Remove unreachable code" example still used the old
#     "goto L2" pseudo-instruction; reword it to "br L2" (as already used
#     on slide 9) by splitting "goto" out of the run and retyping it "br".
$s8 = $p.Slides.Item(8)
$sh8 = $s8.Shapes.Item(3)
$tr8 = $sh8.TextFrame.TextRange
$goto = $tr8.Find("goto", 0)
$goto.Text = "br"

# --- Slide 9 ("Flow control optimization" bullets): the "L1: br L2" line
#     had " " and "L2" split into two separate (but identically formatted)
#     runs; collapse them back into a single " L2" run.
$s9 = $p.Slides.Item(9)
$sh9 = $s9.Shapes.Item(3)
$tr9 = $sh9.TextFrame.TextRange
$spaceL2 = $tr9.Find(" L2", 0)
$spaceL2.Text = " L2"
